# Regenerate save_data to use K instead of Strike# for the "K" column (column G)
# Update the affected rows' K values to the newly computed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 0
$ws.Range("G5").Value = 1
$ws.Range("G7").Value = 0
